# Auto-generated edit script reproducing the crypto-price refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.820.87'
$ws.Range('E2').Value = '  +2.39%  '
$ws.Range('D3').Value = '3.566.92'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''581.35'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').Value = '''187.63'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('D7').Value = '''0.627'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('D8').Value = '3.556.44'
$ws.Range('E8').Value = '  +1.45%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  +17.85%  '
$ws.Range('D11').Value = '''0.652'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').Value = '''54.50'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('E13').Value = '  +5.88%  '
$ws.Range('D14').Value = '''9.55'
$ws.Range('E14').Value = '  +0.93%  '
$ws.Range('D15').Value = '4.135.76'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('D16').Value = '70.859.16'
$ws.Range('E16').Value = '  +2.52%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '''12.86'
$ws.Range('E17').Value = '  +4.30%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '''19.24'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').Value = '3.558.46'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').Value = '''573.26'
$ws.Range('E20').Value = '  +5.35%  '
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').Value = '''1.00'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').Value = '''17.71'
$ws.Range('E23').Value = '  -3.65%  '
$ws.Range('E24').Value = '  +3.39%  '
$ws.Range('D25').Value = '''4.89'
$ws.Range('E25').Value = '  -2.15%  '
$ws.Range('D26').Value = '''94.12'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').Value = '''11.22'
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('E29').Value = '  +2.49%  '
$ws.Range('D30').Value = '''32.80'
$ws.Range('E30').Value = '  +3.18%  '
$ws.Range('D31').Value = '''7.23'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').Value = '''12.33'
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('D33').Value = '''0.116'
$ws.Range('E33').Value = '  +2.28%  '
$ws.Range('E34').Value = '  +23.55%  '
$ws.Range('D35').Value = '''63.23'
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('D36').Value = '''3.30'
$ws.Range('E36').Value = '  +6.00%  '
$ws.Range('D37').Value = '''537.84'
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('E38').Value = '  +2.38%  '
$ws.Range('D39').Value = '0.0₃0817'
$ws.Range('E39').Value = '  +6.66%  '
$ws.Range('D40').Value = '''38.15'
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').Value = '3.631.70'
$ws.Range('E42').Value = '  +10.72%  '
$ws.Range('E43').Value = '  +5.10%  '
$ws.Range('E44').Value = '  +2.37%  '
$ws.Range('E45').Value = '  +6.13%  '
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('D47').Value = '''3.46'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('E48').Value = '  +5.01%  '
$ws.Range('E49').Value = '  +2.79%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = '''1.45'
$ws.Range('E51').Value = '  +3.97%  '
